$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(-0.08312210549351147, 0.6473643844303238, 0.6561623056018049, 0.8100384593349904, 0.8340433064168259, 15)
    3  = @(-0.01572770659491909, 0.4445496780545906, 0.3664663352916562, 0.6053646300302457, 0.6280045102432898, 14)
    4  = @(-0.06939138655364743, 0.5952380793966869, 0.6460436807749044, 0.803768424843191, 0.8334651917952676, 13)
    5  = @(0.06184427096005268, 0.33700112578256, 0.209417853196624, 0.4576219544521701, 0.4735857100263252, 12)
    6  = @(-0.06798794410060216, 0.4968139668552824, 0.377721656252307, 0.6145906412013666, 0.640631895303542, 11)
    7  = @(0.01694414359772521, 0.2607354299704002, 0.1060912636732037, 0.3257165388389169, 0.3428704965612517, 10)
    8  = @(-0.06054918682133375, 0.3330393258437492, 0.1559239451819853, 0.3948720617896198, 0.4138719111052573, 9)
    9  = @(-0.01715367655228583, 0.2997021343775136, 0.1700925356682539, 0.4124227632760514, 0.4405169522264152, 8)
    10 = @(-0.05536040423422233, 0.4661046776547557, 0.3164365824987914, 0.5625269615749909, 0.6046490245027418, 7)
    11 = @(-0.0196809274911882, 0.3184982995745606, 0.1151557242996088, 0.3393460244346599, 0.371109232533003, 6)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
